$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 29043
$ws.Range("E2").Value = 564859265747
$ws.Range("F2").Value = 4686382596
$ws.Range("G2").Value = 0.1864

$ws.Range("D3").Value = 1830.57
$ws.Range("E3").Value = 219940033289
$ws.Range("F3").Value = 2850334802
$ws.Range("G3").Value = 0.04394

$ws.Range("D4").Value = 0.998856
$ws.Range("E4").Value = 83807808883
$ws.Range("F4").Value = 7890074017
$ws.Range("G4").Value = 0.06114

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'BNB'
$ws.Range("D5").Value = 244.17
$ws.Range("E5").Value = 37561392502
$ws.Range("F5").Value = 321361698
$ws.Range("G5").Value = 1.1592

$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'XRP'
$ws.Range("D6").Value = 0.631925
$ws.Range("E6").Value = 33334014662
$ws.Range("F6").Value = 808895231
$ws.Range("G6").Value = 1.10985

$ws.Range("E7").Value = 26051364257
$ws.Range("F7").Value = 1819552024
$ws.Range("G7").Value = 0.05179

$ws.Range("D8").Value = 1828.56
$ws.Range("E8").Value = 14601568803
$ws.Range("F8").Value = 7631986
$ws.Range("G8").Value = 0.04419

$ws.Range("D9").Value = 0.07517699999999999
$ws.Range("E9").Value = 10558746261
$ws.Range("F9").Value = 443347385
$ws.Range("G9").Value = -1.24

$ws.Range("D10").Value = 0.294552
$ws.Range("E10").Value = 10321851414
$ws.Range("F10").Value = 158873178
$ws.Range("G10").Value = 1.39015

$ws.Range("D11").Value = 23.2
$ws.Range("E11").Value = 9401949422
$ws.Range("F11").Value = 309850817
$ws.Range("G11").Value = 2.26719

$ws.Range("D12").Value = 0.077129
$ws.Range("E12").Value = 6906643033
$ws.Range("F12").Value = 163162283
$ws.Range("G12").Value = -0.35102

$ws.Range("B13").Value = 'DOT'
$ws.Range("C13").Value = 'Polkadot'
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = 6311831344
$ws.Range("F13").Value = 71916117
$ws.Range("G13").Value = 1.1746

$ws.Range("D14").Value = 0.668238
$ws.Range("E14").Value = 6228832397
$ws.Range("F14").Value = 174531251
$ws.Range("G14").Value = 0.92513

$ws.Range("B15").Value = 'LTC'
$ws.Range("C15").Value = 'Litecoin'
$ws.Range("D15").Value = 83.26000000000001
$ws.Range("E15").Value = 6120130955
$ws.Range("F15").Value = 282178143
$ws.Range("G15").Value = 1.19595

$ws.Range("D16").Value = [double]"9.77e-06"
$ws.Range("E16").Value = 5770889760
$ws.Range("F16").Value = 833295208
$ws.Range("G16").Value = 3.36204

$ws.Range("B17").Value = 'WBTC'
$ws.Range("C17").Value = 'Wrapped Bitcoin'
$ws.Range("D17").Value = 29037
$ws.Range("E17").Value = 4710158683
$ws.Range("F17").Value = 30287320
$ws.Range("G17").Value = 0.24604

$ws.Range("B18").Value = 'UNI'
$ws.Range("C18").Value = 'Uniswap'
$ws.Range("D18").Value = 6.06
$ws.Range("E18").Value = 4562970588
$ws.Range("F18").Value = 60049745
$ws.Range("G18").Value = 1.5731

$ws.Range("D19").Value = 225.82
$ws.Range("E19").Value = 4395929329
$ws.Range("F19").Value = 130978863
$ws.Range("G19").Value = 0.93953

$ws.Range("D20").Value = 12.6
$ws.Range("E20").Value = 4363320291
$ws.Range("F20").Value = 92117487
$ws.Range("G20").Value = 2.46253

$ws.Range("B21").Value = 'TON'
$ws.Range("C21").Value = 'Toncoin'
$ws.Range("D21").Value = 1.2
$ws.Range("E21").Value = 4141467095
$ws.Range("F21").Value = 21567840
$ws.Range("G21").Value = 0.49407

$ws.Range("B22").Value = 'DAI'
$ws.Range("C22").Value = 'Dai'
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 4044572604
$ws.Range("F22").Value = 68618576
$ws.Range("G22").Value = 0.10263

$ws.Range("B23").Value = 'XLM'
$ws.Range("C23").Value = 'Stellar'
$ws.Range("D23").Value = 0.141593
$ws.Range("E23").Value = 3869370199
$ws.Range("F23").Value = 107226817
$ws.Range("G23").Value = 4.26286

$ws.Range("B24").Value = 'LINK'
$ws.Range("C24").Value = 'Chainlink'
$ws.Range("D24").Value = 7.14
$ws.Range("E24").Value = 3843051482
$ws.Range("F24").Value = 131028854
$ws.Range("G24").Value = -0.45324

$ws.Range("D25").Value = 4.02
$ws.Range("E25").Value = 3731907083
$ws.Range("F25").Value = 1825781
$ws.Range("G25").Value = -0.19076

$ws.Range("D26").Value = 0.9997239999999999
$ws.Range("E26").Value = 3415836319
$ws.Range("F26").Value = 983733960
$ws.Range("G26").Value = 0.0266

$ws.Range("D27").Value = 0.999098
$ws.Range("E27").Value = 2962479271
$ws.Range("F27").Value = 852152725
$ws.Range("G27").Value = -0.05358

$ws.Range("D28").Value = 160.31
$ws.Range("E28").Value = 2917536909
$ws.Range("F28").Value = 89069659
$ws.Range("G28").Value = 0.43787

$ws.Range("B29").Value = 'OKB'
$ws.Range("C29").Value = 'OKB'
$ws.Range("D29").Value = 45.55
$ws.Range("E29").Value = 2732843518
$ws.Range("F29").Value = 5467927
$ws.Range("G29").Value = 0.35866

$ws.Range("D30").Value = 17.94
$ws.Range("E30").Value = 2554069026
$ws.Range("F30").Value = 52734198
$ws.Range("G30").Value = 0.93449

$ws.Range("B31").Value = 'ATOM'
$ws.Range("C31").Value = 'Cosmos Hub'
$ws.Range("D31").Value = 8.52
$ws.Range("E31").Value = 2491798732
$ws.Range("F31").Value = 103870626
$ws.Range("G31").Value = 1.43999

$ws.Range("D32").Value = 4.13
$ws.Range("E32").Value = 1814013626
$ws.Range("F32").Value = 59806171
$ws.Range("G32").Value = 1.99424

$ws.Range("D33").Value = 4.06
$ws.Range("E33").Value = 1790335119
$ws.Range("F33").Value = 13610391
$ws.Range("G33").Value = 1.14665

$ws.Range("D34").Value = 0.054888
$ws.Range("E34").Value = 1787889471
$ws.Range("F34").Value = 29697902
$ws.Range("G34").Value = 5.70609

$ws.Range("B35").Value = 'MNT'
$ws.Range("C35").Value = 'Mantle'
$ws.Range("D35").Value = 0.508229
$ws.Range("E35").Value = 1644115400
$ws.Range("F35").Value = 7948757
$ws.Range("G35").Value = -0.54503

$ws.Range("B36").Value = 'LDO'
$ws.Range("C36").Value = 'Lido DAO'
$ws.Range("D36").Value = 1.86
$ws.Range("E36").Value = 1634824146
$ws.Range("F36").Value = 54104342
$ws.Range("G36").Value = 1.10114

$ws.Range("B37").Value = 'CRO'
$ws.Range("C37").Value = 'Cronos'
$ws.Range("D37").Value = 0.058012
$ws.Range("E37").Value = 1521078722
$ws.Range("F37").Value = 4525889
$ws.Range("G37").Value = 1.09347

$ws.Range("B38").Value = 'APT'
$ws.Range("C38").Value = 'Aptos'
$ws.Range("D38").Value = 6.78
$ws.Range("E38").Value = 1486794977
$ws.Range("F38").Value = 42734260
$ws.Range("G38").Value = 1.46874

$ws.Range("B39").Value = 'QNT'
$ws.Range("C39").Value = 'Quant'
$ws.Range("D39").Value = 101.88
$ws.Range("E39").Value = 1481368981
$ws.Range("F39").Value = 16212791
$ws.Range("G39").Value = 0.48671

$ws.Range("D40").Value = 1.14
$ws.Range("E40").Value = 1451998905
$ws.Range("F40").Value = 114484250
$ws.Range("G40").Value = -0.57778

$ws.Range("D41").Value = 0.01785485
$ws.Range("E41").Value = 1296451013
$ws.Range("F41").Value = 25984554
$ws.Range("G41").Value = 0.26489

$ws.Range("D42").Value = 1.36
$ws.Range("E42").Value = 1277768501
$ws.Range("F42").Value = 37264658
$ws.Range("G42").Value = 1.26672

$ws.Range("B43").Value = 'OP'
$ws.Range("C43").Value = 'Optimism'
$ws.Range("D43").Value = 1.69
$ws.Range("E43").Value = 1214847082
$ws.Range("F43").Value = 91062255
$ws.Range("G43").Value = 2.10258

$ws.Range("B44").Value = 'MKR'
$ws.Range("C44").Value = 'Maker'
$ws.Range("D44").Value = 1246.18
$ws.Range("E44").Value = 1123296555
$ws.Range("F44").Value = 150814405
$ws.Range("G44").Value = -1.48139

$ws.Range("B45").Value = 'KAS'
$ws.Range("C45").Value = 'Kaspa'
$ws.Range("D45").Value = 0.051558
$ws.Range("E45").Value = 1035295811
$ws.Range("F45").Value = 26209171
$ws.Range("G45").Value = 12.4698

$ws.Range("B46").Value = 'XDC'
$ws.Range("C46").Value = 'XDC Network'
$ws.Range("D46").Value = 0.07346
$ws.Range("E46").Value = 1020134211
$ws.Range("F46").Value = 23573112
$ws.Range("G46").Value = 1.75388

$ws.Range("D47").Value = 1977.08
$ws.Range("E47").Value = 978869063
$ws.Range("F47").Value = 2153055
$ws.Range("G47").Value = 0.1216

$ws.Range("B48").Value = 'GRT'
$ws.Range("C48").Value = 'The Graph'
$ws.Range("D48").Value = 0.106285
$ws.Range("E48").Value = 970398056
$ws.Range("F48").Value = 58437186
$ws.Range("G48").Value = 1.97488

$ws.Range("B49").Value = 'AAVE'
$ws.Range("C49").Value = 'Aave'
$ws.Range("D49").Value = 65.40000000000001
$ws.Range("E49").Value = 949423102
$ws.Range("F49").Value = 61364375
$ws.Range("G49").Value = 1.71488

$ws.Range("D50").Value = 0.107631
$ws.Range("E50").Value = 840906547
$ws.Range("F50").Value = 20316973
$ws.Range("G50").Value = 2.24045

$ws.Range("B51").Value = 'SAND'
$ws.Range("C51").Value = 'The Sandbox'
$ws.Range("D51").Value = 0.406102
$ws.Range("E51").Value = 834622081
$ws.Range("F51").Value = 70818243
$ws.Range("G51").Value = 2.33244

